# Update the cryptocurrency price/volume snapshot on Sheet1.
# Column D (Price) values are stored as plain text in this sheet, so any
# new value that looks like a number is entered with a leading apostrophe
# to keep Excel from converting it to a numeric cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.966.13"
$ws.Range("E2").Value = "  -5.94%  "

$ws.Range("D3").Value = "3.247.60"
$ws.Range("E3").Value = "  -7.29%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'177.10"
$ws.Range("E5").Value = "  -12.08%  "

$ws.Range("D6").Value = "'516.12"
$ws.Range("E6").Value = "  -6.62%  "

$ws.Range("E7").Value = "  -1.06%  "

$ws.Range("D8").Value = "3.239.84"
$ws.Range("E8").Value = "  -7.19%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "'0.615"
$ws.Range("E10").Value = "  -6.11%  "

$ws.Range("D11").Value = "'57.90"
$ws.Range("E11").Value = "  -4.65%  "

$ws.Range("E12").Value = "  -8.59%  "

$ws.Range("D13").Value = "'0.0000254"
$ws.Range("E13").Value = "  -7.02%  "

$ws.Range("D14").Value = "'9.06"
$ws.Range("E14").Value = "  -7.36%  "

$ws.Range("D15").Value = "3.770.31"
$ws.Range("E15").Value = "  -7.35%  "

$ws.Range("D16").Value = "'0.116"
$ws.Range("E16").Value = "  -6.26%  "

$ws.Range("D17").Value = "3.252.57"
$ws.Range("E17").Value = "  -7.17%  "

$ws.Range("D18").Value = "'17.46"
$ws.Range("E18").Value = "  -5.04%  "

$ws.Range("D19").Value = "62.873.73"
$ws.Range("E19").Value = "  -5.71%  "

$ws.Range("D20").Value = "'10.88"
$ws.Range("E20").Value = "  -7.63%  "

$ws.Range("D21").Value = "'0.946"
$ws.Range("E21").Value = "  -8.37%  "

$ws.Range("D22").Value = "'369.74"
$ws.Range("E22").Value = "  -4.74%  "

$ws.Range("D23").Value = "'11.15"
$ws.Range("E23").Value = "  -6.49%  "

$ws.Range("E24").Value = "  -7.32%  "

$ws.Range("D25").Value = "'79.80"
$ws.Range("E25").Value = "  -3.01%  "

$ws.Range("D26").Value = "'3.82"
$ws.Range("E26").Value = "  +3.48%  "

$ws.Range("E27").Value = "  -1.21%  "

$ws.Range("D28").Value = "'11.36"
$ws.Range("E28").Value = "  -4.41%  "

$ws.Range("E29").Value = "  -6.04%  "

$ws.Range("D30").Value = "'8.27"
$ws.Range("E30").Value = "  -6.30%  "

$ws.Range("D31").Value = "'28.40"
$ws.Range("E31").Value = "  -6.99%  "

$ws.Range("D32").Value = "'6.82"
$ws.Range("E32").Value = "  -6.57%  "

$ws.Range("D33").Value = "'631.26"
$ws.Range("E33").Value = "  -5.98%  "

$ws.Range("D34").Value = "'11.25"
$ws.Range("E34").Value = "  -3.73%  "

$ws.Range("D35").Value = "'0.105"
$ws.Range("E35").Value = "  -3.72%  "

$ws.Range("D36").Value = "'58.30"
$ws.Range("E36").Value = "  -7.84%  "

$ws.Range("D37").Value = "'0.401"
$ws.Range("E37").Value = "  -1.49%  "

$ws.Range("E38").Value = "  +0.16%  "

$ws.Range("D39").Value = "'36.34"
$ws.Range("E39").Value = "  -7.10%  "

$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  -0.04%  "

$ws.Range("D41").Value = "2.933.99"
$ws.Range("E41").Value = "  -6.52%  "

$ws.Range("D42").Value = "'0.124"
$ws.Range("E42").Value = "  -3.54%  "

$ws.Range("D43").Value = "0.0₃0660"
$ws.Range("E43").Value = "  -6.15%  "

$ws.Range("D44").Value = "'2.44"
$ws.Range("E44").Value = "  -3.52%  "

$ws.Range("D45").Value = "'2.67"
$ws.Range("E45").Value = "  -12.65%  "

$ws.Range("D46").Value = "'0.0391"
$ws.Range("E46").Value = "  -1.46%  "

$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").Value = "'2.82"
$ws.Range("E47").Value = "  +8.44%  "

$ws.Range("B48").Value = "WEMIXToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "'2.58"
$ws.Range("E48").Value = "  -4.29%  "

$ws.Range("B49").Value = "ApeXProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D49").Value = "'2.97"
$ws.Range("E49").Value = "  +2.30%  "

$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D50").Value = "'0.125"
$ws.Range("E50").Value = "  -1.55%  "

$ws.Range("E51").Value = "  -11.07%  "
